$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. A5: "John Deery" -> "Branan Harrison" (name correction, style unchanged)
$ws.Range("A5").Value = "Branan Harrison"

# 2. C8: "9:45AM-4PM" -> "OFF" (style goes from 4 (no fill) to 6 (OFF yellow fill))
#    Copy the format from a neighboring OFF cell (B8, style 6) then set the text.
$ws.Range("B8").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C8").Value = "OFF"

# 3. C12: "9:30AM-4PM" -> "9:45AM-4PM" (style 4 stays the same, text only)
$ws.Range("C12").Value = "9:45AM-4PM"

# 4. C15: "OFF" -> "10AM-5PM" (style goes from 6 (OFF yellow fill) to 4 (no fill))
#    Copy the format from a neighboring working-shift cell (B15, style 4) then set the text.
$ws.Range("B15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C15").Value = "10AM-5PM"

# 5. C16: "10AM-5PM" -> "OFF" (style goes from 4 (no fill) to 6 (OFF yellow fill))
#    Copy the format from a neighboring OFF cell (B16, style 6) then set the text.
$ws.Range("B16").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C16").Value = "OFF"

$excel.CutCopyMode = 0

# 6. Unassigned-shifts block (rows 27-29): C27 changes from "Bartender, 10AM-4PM"
#    to "Lifeguard, 9:30AM-4PM", and the old C27 content is now placed at the
#    previously-empty C28 (format copied from the already-styled B28 cell).
$ws.Range("B28").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C28").Value = "Bartender,`n10AM-4PM"
$excel.CutCopyMode = 0

$ws.Range("C27").Value = "Lifeguard,`n9:30AM-4PM"
